$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44; rows 44-161 shift down to 45-162.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new record.
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = 45251
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100112010
$ws.Range("G44").Value = "Achicoria"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 90
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("M44").Value = 10000
$ws.Range("N44").Value = "$/caja 18 unidades"
$ws.Range("O44").Value = "Región Metropolitana"
$ws.Range("P44").Value = 556
$ws.Range("Q44").Value = 18
$ws.Range("R44").Value = "Hortaliza"

# Preserve the date-number style for the new date cell (column D uses style index 2).
$ws.Range("D44").NumberFormat = $ws.Range("D45").NumberFormat
